# risorse.xlsx — "Add files via upload" edit
#
# The commit inserts two new resource rows ("B212" and "B217") into the
# existing risorsa/dimensione table on Sheet1, right-aligns the
# "dimensione" (column B) values, and updates the current on-screen
# selection to the bottom of the refreshed table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert "B212" (dimensione 4) right after "B211" (old row 161) ---
# Before the edit, row 162 held "B213"; pushing it (and everything below)
# down by one row makes room for the new "B212" entry.
$ws.Rows("162:162").Insert()
$ws.Range("A162").Value = "B212"
$ws.Range("B162").Value = 4
$ws.Range("B162").NumberFormat = $ws.Range("B161").NumberFormat
$ws.Range("B162").HorizontalAlignment = $ws.Range("B161").HorizontalAlignment

# --- 2. Insert "B217" (dimensione 4) right after "B215" ---
# After step 1, the old "B218" row now lives at row 166; insert above it.
$ws.Rows("166:166").Insert()
$ws.Range("A166").Value = "B217"
$ws.Range("B166").Value = 4
$ws.Range("B166").NumberFormat = $ws.Range("B165").NumberFormat
$ws.Range("B166").HorizontalAlignment = $ws.Range("B165").HorizontalAlignment

# --- 3. Right-align the whole "dimensione" data column (was left-aligned) ---
# This covers every row that uses the same number format as the two rows
# just inserted, including them, from the first "style 3" row through the
# new last row of the (now 195-row) table.
$ws.Range("B152:B195").HorizontalAlignment = -4152

# --- 4. Refresh the active selection/view to sit on the new last data row ---
$ws.Range("A190").Select()
$excel.ActiveWindow.ScrollRow = 180
$excel.ActiveWindow.ScrollColumn = 1
